$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "[sound:0003_ഒപ്പം_02.mp3]"
$ws.Range("G4").Value = "[sound:0003_ഒപ്പം_03.mp3]"
$ws.Range("G5").Value = "[sound:0003_ഒപ്പം_04.mp3]"
$ws.Range("G6").Value = "[sound:0003_ഒപ്പം_05.mp3]"
$ws.Range("G7").Value = "[sound:0003_ഒപ്പം_06.mp3]"
$ws.Range("G8").Value = "[sound:0003_ഒപ്പം_07.mp3]"
$ws.Range("G9").Value = "[sound:0003_ഒപ്പം_08.mp3]"
$ws.Range("G10").Value = "[sound:0003_ഒപ്പം_09.mp3]"
$ws.Range("G11").Value = "[sound:0003_ഒപ്പം_10.mp3]"
